# Update the "Add Devices" test data sheet: change the Node/Panel Node
# values in B8:B11 from "Node1" to "Node" and move the selection to B9:B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")
$ws.Activate()

$ws.Range("B8:B11").Value = "Node"

$ws.Range("B9:B11").Select()
